$wb = $excel.ActiveWorkbook

# --- SEC_Processes: rename the hard-coal mine process to wind mine -----
$wsProc = $wb.Worksheets.Item("SEC_Processes")
$wsProc.Range("D9").Value = "MIN_EX_WIND_ON"
$wsProc.Range("E9").Value = "Wind Mine"

# --- SEC_Comm: rename the HARD_COAL commodity to WIND_ON ---------------
$wsComm = $wb.Worksheets.Item("SEC_Comm")
$wsComm.Range("C9").Value = "WIND_ON"
$wsComm.Range("D9").Value = "Wind Onshore"

# --- SEC_Processes: remove the Kozienice power-plant process row -------
$wsProc.Range("B10").Value = $null
$wsProc.Range("D10").Value = $null
$wsProc.Range("E10").Value = $null
$wsProc.Range("F10").Value = $null
$wsProc.Range("G10").Value = $null
$wsProc.Range("H10").Value = $null

# --- MIN_IMP: update extraction cost / bound for the wind mine ---------
$wsMin = $wb.Worksheets.Item("MIN_IMP")
$wsMin.Range("E9").Value = 0.001
$wsMin.Range("E9").NumberFormat = "0.000"
$wsMin.Range("F9").Value = $null

# --- PP: remove the Kozienice power-plant transformation row -----------
$wsPP = $wb.Worksheets.Item("PP")
$wsPP.Range("B9").Value = $null
$wsPP.Range("C9").Value = $null
$wsPP.Range("D9").Value = $null
$wsPP.Range("E9").Value = $null
$wsPP.Range("F9").Value = $null
$wsPP.Range("G9").Value = $null
$wsPP.Range("H9").Value = $null
$wsPP.Range("I9").Value = $null
$wsPP.Range("J9").Value = $null
$wsPP.Range("K9").Value = $null
